$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("H2").Value = [double]"1.315819881037222e-16"
$ws.Range("K2").Value = 58.20563568893839
$ws.Range("L2").Value = "[53.87264349144694, 62.53862788642985]"
$ws.Range("O2").Value = 1.515763422452733
$ws.Range("P2").Value = "[1.4402897250691948, 1.5912371198362711]"
$ws.Range("S2").Value = 52.50961682556846
$ws.Range("T2").Value = "[49.63289466093811, 55.386338990198816]"
$ws.Range("W2").Value = 17.58044044044058
$ws.Range("X2").Value = 17.30212212212226
$ws.Range("Y2").Value = 17.8587587587589

# Row 3 updates
$ws.Range("E3").Value = 22.6300000000001
$ws.Range("H3").Value = [double]"1.315819881037222e-16"
$ws.Range("K3").Value = 57.00602027131258
$ws.Range("L3").Value = "[50.570708689039655, 63.4413318535855]"
$ws.Range("O3").Value = 0.3207632138800394
$ws.Range("P3").Value = "[0.19497371824080822, 0.4465527095192705]"
$ws.Range("Q3").Value = [double]"1.135338931090857e-06"
$ws.Range("R3").Value = [double]"1.135338931090857e-06"
$ws.Range("S3").Value = 52.84706368075487
$ws.Range("T3").Value = "[49.19274803266827, 56.50137932884147]"
$ws.Range("W3").Value = 21.47471471471481
$ws.Range("X3").Value = 21.02166166166175
$ws.Range("Y3").Value = 21.92776776776786
